$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

# Columns B through K on row 23 were empty inline strings; set them to the
# text "nan" to match the pattern used in the rows above (18-22).
$cols = @("B","C","D","E","F","G","H","I","J","K")
foreach ($col in $cols) {
    $ws.Range($col + "23").Value = "nan"
}

# Update the correction text in N23, changing (988) to (993.6).
$ws.Range("N23").Value = "(993.6)تغيير الفلاتس المتحركه و جريده الخلفيه رقم (1) عند طن"
